# Update countries & provincias Spain
# Applies the 2020-04-14 17:52 COVID-19 snapshot refresh:
#  - updates the "Casos totales" A1 title timestamp
#  - refreshes per-country stats for the countries whose totals changed
#  - re-sorts all data rows (4-216) by "Casos totales" (column B) descending,
#    matching the country order produced by the refreshed dataset

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full refreshed table in final (sorted) row order: Pais, Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$data = @(
    @('Estados Unidos', 589015, 2074, 37326, 527162, 12772, 887, 24527),
    @('España', 172541, 2442, 67504, 86981, 7371, 300, 18056),
    @('Italia', 159516, 0, 35435, 103616, 3260, 0, 20465),
    @('Francia', 136779, 0, 27718, 94094, 6821, 0, 14967),
    @('Alemania', 130694, 622, 68200, 59233, 4288, 67, 3261),
    @('Reino Unido', 93873, 5252, 0, 81422, 1559, 778, 12107),
    @('China', 82249, 89, 77738, 1170, 116, 0, 3341),
    @('Iran', 74877, 1574, 48129, 22065, 3691, 98, 4683),
    @('Turquia', 61049, 0, 3957, 55796, 1786, 0, 1296),
    @('Belgica', 31119, 530, 6868, 20094, 1223, 254, 4157),
    @('Paises Bajos', 27419, 868, 250, 24224, 1358, 122, 2945),
    @('Canada', 26163, 483, 7756, 17584, 557, 43, 823),
    @('Suiza', 25913, 225, 13700, 11051, 386, 24, 1162),
    @('Brasil', 23955, 525, 2979, 19615, 296, 33, 1361),
    @('Rusia', 21102, 2774, 1694, 19238, 8, 22, 170),
    @('Portugal', 17448, 514, 347, 16534, 218, 32, 567),
    @('Austria', 14185, 144, 7633, 6168, 243, 16, 384),
    @('Israel', 11868, 282, 2000, 9751, 181, 1, 117),
    @('Suecia', 11445, 497, 381, 10031, 915, 114, 1033),
    @('Irlanda', 10647, 0, 25, 10257, 194, 0, 365),
    @('Corea del Sur', 10564, 27, 7534, 2808, 55, 5, 222),
    @('India', 10541, 88, 1205, 8978, 0, 0, 358),
    @('Peru', 9784, 0, 2642, 6926, 143, 0, 216),
    @('Chile', 7917, 392, 2646, 5179, 387, 10, 92),
    @('Japon', 7645, 27, 799, 6703, 135, 0, 143),
    @('Ecuador', 7529, 0, 597, 6577, 121, 0, 355),
    @('Polonia', 7202, 268, 618, 6321, 160, 18, 263),
    @('Rumania', 6879, 246, 1051, 5477, 241, 20, 351),
    @('Noruega', 6623, 20, 32, 6452, 59, 5, 139),
    @('Dinamarca', 6496, 178, 2515, 3682, 100, 14, 299),
    @('Australia', 6400, 41, 3598, 2741, 80, 0, 61),
    @('Chequia', 6101, 42, 642, 5298, 92, 18, 161),
    @('Pakistan', 5837, 341, 1378, 4363, 46, 3, 96),
    @('Arabia Saudita', 5369, 435, 889, 4407, 59, 8, 73),
    @('Filipinas', 5223, 291, 295, 4593, 1, 20, 335),
    @('Mexico', 5014, 353, 1964, 2718, 207, 36, 332),
    @('Malasia', 4987, 170, 2478, 2427, 60, 5, 82),
    @('Indonesia', 4839, 282, 426, 3954, 0, 60, 459),
    @('Emiratos Arabes Unidos', 4521, 0, 852, 3644, 1, 0, 25),
    @('Serbia', 4465, 411, 400, 3971, 131, 9, 94),
    @('Panama', 3472, 0, 61, 3317, 106, 0, 94),
    @('Catar', 3428, 197, 373, 3048, 37, 0, 7),
    @('Ucrania', 3372, 270, 119, 3155, 45, 5, 98),
    @('Luxemburgo', 3292, 0, 500, 2723, 30, 0, 69),
    @('Republica Dominicana', 3286, 119, 162, 2941, 143, 6, 183),
    @('Bielorrusia', 3281, 362, 203, 3045, 57, 4, 33),
    @('Singapur', 3252, 334, 611, 2631, 29, 1, 10),
    @('Finlandia', 3161, 97, 300, 2797, 74, 5, 64),
    @('Colombia', 2852, 0, 319, 2421, 106, 0, 112),
    @('Tailandia', 2613, 34, 1405, 1167, 61, 1, 41),
    @('Sudafrica', 2415, 143, 410, 1978, 7, 0, 27),
    @('Argentina', 2277, 0, 559, 1617, 83, 3, 101),
    @('Egipto', 2190, 0, 589, 1437, 0, 0, 164),
    @('Grecia', 2170, 25, 269, 1800, 76, 2, 101),
    @('Argelia', 1983, 0, 601, 1069, 60, 0, 313),
    @('Moldavia', 1934, 222, 134, 1764, 80, 1, 36),
    @('Marruecos', 1838, 75, 210, 1502, 1, 0, 126),
    @('Islandia', 1720, 9, 989, 723, 8, 0, 8),
    @('Croacia', 1704, 54, 415, 1258, 34, 6, 31),
    @('Barein', 1522, 161, 645, 870, 3, 1, 7),
    @('Hungria', 1512, 54, 122, 1268, 58, 13, 122),
    @('Irak', 1400, 22, 766, 556, 0, 0, 78),
    @('Estonia', 1373, 41, 115, 1227, 11, 3, 31),
    @('Nueva Zelanda', 1366, 17, 628, 729, 4, 4, 9),
    @('Kuwait', 1355, 55, 176, 1176, 26, 1, 3),
    @('Eslovenia', 1220, 8, 152, 1012, 35, 1, 56),
    @('Kazajistan', 1202, 111, 203, 985, 21, 2, 14),
    @('Azerbaiyan', 1197, 49, 351, 833, 25, 1, 13),
    @('Uzbekistan', 1113, 115, 99, 1010, 8, 0, 4),
    @('Bosnia y Herzegovina', 1083, 46, 236, 807, 4, 1, 40),
    @('Lituania', 1070, 8, 101, 943, 14, 2, 26),
    @('Armenia', 1067, 28, 265, 786, 30, 2, 16),
    @('Hong Kong', 1013, 3, 434, 575, 13, 0, 4),
    @('Banglades', 1012, 209, 42, 924, 1, 7, 46),
    @('Republica de Macedonia', 908, 54, 86, 778, 15, 6, 44),
    @('Camerun', 848, 0, 130, 704, 0, 2, 14),
    @('Eslovaquia', 835, 66, 113, 720, 5, 0, 2),
    @('Oman', 813, 86, 130, 679, 3, 0, 4),
    @('Cuba', 766, 40, 132, 613, 11, 0, 21),
    @('Tunez', 726, 0, 43, 649, 89, 0, 34),
    @('Afganistan', 714, 49, 40, 651, 0, 2, 23),
    @('Bulgaria', 713, 28, 81, 597, 29, 3, 35),
    @('Crucero', 712, 0, 639, 61, 7, 0, 12),
    @('Republica de Chipre', 695, 33, 65, 618, 8, 0, 12),
    @('Letonia', 657, 2, 16, 636, 3, 0, 5),
    @('Principado de Andorra', 646, 0, 128, 489, 17, 0, 29),
    @('Libano', 641, 9, 80, 540, 33, 1, 21),
    @('Ghana', 636, 70, 17, 611, 2, 0, 8),
    @('Costa de Marfil', 626, 0, 89, 531, 0, 0, 6),
    @('Costa Rica', 612, 0, 62, 547, 14, 0, 3),
    @('Niger', 548, 0, 86, 449, 0, 0, 13),
    @('Burkina Faso', 515, 0, 161, 327, 0, 0, 27),
    @('Uruguay', 483, 0, 248, 227, 15, 0, 8),
    @('Albania', 475, 8, 248, 203, 5, 1, 24),
    @('Kirguistan', 430, 11, 71, 354, 5, 0, 5),
    @('Honduras', 407, 10, 7, 374, 10, 1, 26),
    @('Malta', 393, 9, 44, 346, 4, 0, 3),
    @('Taiwan', 393, 0, 124, 263, 0, 0, 6),
    @('Reunion', 391, 0, 40, 351, 3, 0, 0),
    @('Jordania', 391, 0, 215, 169, 5, 0, 7),
    @('San Marino', 371, 0, 53, 282, 15, 0, 36),
    @('Guinea', 363, 44, 31, 332, 0, 0, 0),
    @('Republica de Yibuti', 363, 65, 53, 308, 0, 0, 2),
    @('Bolivia', 354, 24, 6, 320, 3, 1, 28),
    @('Nigeria', 343, 0, 91, 242, 2, 0, 10),
    @('Mauricio', 324, 0, 51, 264, 3, 0, 9),
    @('Estado de Palestina', 308, 0, 62, 244, 0, 0, 2),
    @('Senegal', 299, 8, 183, 114, 1, 0, 2),
    @('Georgia', 296, 24, 68, 225, 6, 0, 3),
    @('Montenegro', 278, 4, 5, 269, 7, 1, 4),
    @('Vietnam', 266, 1, 169, 97, 8, 0, 0),
    @('Isla de Man', 242, 0, 141, 99, 13, 0, 2),
    @('Consejo Danes para los Refugiados', 241, 6, 20, 201, 0, 0, 20),
    @('Sri Lanka', 219, 2, 61, 151, 1, 0, 7),
    @('Mayotte', 217, 10, 69, 145, 3, 0, 3),
    @('Kenia', 208, 0, 40, 159, 2, 0, 9),
    @('Venezuela', 189, 0, 110, 70, 6, 0, 9),
    @('Islas Feroe', 184, 0, 163, 21, 0, 0, 0),
    @('Guatemala', 167, 11, 19, 143, 3, 0, 5),
    @('Paraguay', 159, 12, 22, 130, 1, 1, 7),
    @('Martinica', 157, 0, 50, 101, 19, 0, 6),
    @('El Salvador', 149, 12, 25, 118, 2, 0, 6),
    @('Guadalupe', 145, 0, 67, 70, 13, 0, 8),
    @('Brunei', 136, 0, 107, 28, 2, 0, 1),
    @('Gibraltar', 129, 0, 93, 36, 1, 0, 0),
    @('Ruanda', 127, 0, 42, 85, 0, 0, 0),
    @('Mali', 123, 0, 26, 87, 0, 0, 10),
    @('Camboya', 122, 0, 91, 31, 1, 0, 0),
    @('Trinidad yTobago', 113, 0, 17, 88, 0, 0, 8),
    @('Madagascar', 108, 2, 23, 85, 1, 0, 0),
    @('Monaco', 93, 0, 6, 86, 5, 0, 1),
    @('Aruba', 92, 0, 32, 60, 1, 0, 0),
    @('Guayana Francesa', 86, 0, 51, 35, 1, 0, 0),
    @('Etiopia', 82, 8, 14, 65, 0, 0, 3),
    @('Liechtenstein', 79, 0, 55, 23, 0, 0, 1),
    @('Togo', 77, 0, 32, 42, 0, 0, 3),
    @('Jamaica', 73, 0, 19, 50, 0, 0, 4),
    @('Barbados', 72, 0, 13, 55, 4, 0, 4),
    @('Birmania', 63, 1, 2, 57, 0, 0, 4),
    @('Somalia', 60, 0, 2, 56, 2, 0, 2),
    @('Congo', 60, 0, 5, 50, 0, 0, 5),
    @('Liberia', 59, 0, 4, 49, 0, 0, 6),
    @('Gabon', 57, 0, 1, 55, 0, 0, 1),
    @('Bermudas', 57, 0, 30, 22, 3, 0, 5),
    @('Polinesia Francesa', 55, 0, 0, 55, 1, 0, 0),
    @('Islas Caimanes', 54, 0, 6, 47, 3, 0, 1),
    @('Uganda', 54, 0, 8, 46, 0, 0, 0),
    @('San Martin (Parte Holandesa)', 52, 2, 5, 38, 2, 0, 9),
    @('Tanzania', 49, 0, 7, 39, 0, 0, 3),
    @('Bahamas', 49, 0, 6, 35, 1, 0, 8),
    @('Guyana', 47, 0, 8, 33, 5, 0, 6),
    @('Macao', 45, 0, 10, 35, 1, 0, 0),
    @('Zambia', 45, 0, 30, 13, 1, 0, 2),
    @('Guinea Ecuatorial', 41, 20, 4, 37, 0, 0, 0),
    @('Haiti', 40, 0, 0, 37, 0, 0, 3),
    @('Puerto Rico', 39, 0, 1, 36, 0, 0, 2),
    @('Guinea-Bisau', 38, 0, 0, 38, 0, 0, 0),
    @('Benin', 35, 0, 18, 16, 0, 0, 1),
    @('Eritrea', 34, 0, 0, 34, 0, 0, 0),
    @('Guam', 32, 0, 0, 31, 0, 0, 1),
    @('San Martin (Parte Francesa)', 32, 0, 11, 19, 5, 0, 2),
    @('Mongolia', 30, 13, 5, 25, 0, 0, 0),
    @('Siria', 29, 4, 5, 22, 0, 0, 2),
    @('Sudan', 29, 0, 4, 21, 0, 0, 4),
    @('Mozambique', 28, 7, 2, 26, 0, 0, 0),
    @('Libia', 26, 0, 9, 16, 0, 0, 1),
    @('Republica del Chad', 23, 0, 2, 21, 0, 0, 0),
    @('Antigua y Barbuda', 23, 0, 3, 18, 1, 0, 2),
    @('Maldivas', 20, 0, 14, 6, 0, 0, 0),
    @('Laos', 19, 0, 1, 18, 0, 0, 0),
    @('Angola', 19, 0, 4, 13, 0, 0, 2),
    @('Nueva Caledonia', 18, 0, 1, 17, 0, 0, 0),
    @('Belice', 18, 0, 0, 16, 1, 0, 2),
    @('Islas Virgenes de los Estados Unidos', 17, 0, 0, 17, 0, 0, 0),
    @('Zimbabue', 17, 0, 0, 14, 0, 0, 3),
    @('Fiyi', 16, 0, 0, 16, 0, 0, 0),
    @('Nepal', 16, 2, 1, 15, 0, 0, 0),
    @('Malaui', 16, 0, 0, 14, 1, 0, 2),
    @('Namibia', 16, 0, 3, 13, 0, 0, 0),
    @('Dominica', 16, 0, 8, 8, 0, 0, 0),
    @('Santa Lucia', 15, 0, 4, 11, 0, 0, 0),
    @('Suazilandia', 15, 0, 8, 7, 0, 0, 0),
    @('Granada', 14, 0, 0, 14, 2, 0, 0),
    @('Curazao', 14, 0, 10, 3, 0, 0, 1),
    @('Botsuana', 13, 0, 0, 12, 0, 0, 1),
    @('San Cristobal y Nieves', 12, 0, 0, 12, 0, 0, 0),
    @('San Vicente y las Granadinas', 12, 0, 1, 11, 0, 0, 0),
    @('Seychelles', 11, 0, 0, 11, 0, 0, 0),
    @('Sierra Leona', 11, 1, 0, 11, 0, 0, 0),
    @('Islas Malvinas', 11, 6, 1, 10, 0, 0, 0),
    @('Montserrat', 11, 0, 1, 10, 1, 0, 0),
    @('Republica de Africa Central', 11, 0, 4, 7, 0, 0, 0),
    @('Groenlandia', 11, 0, 11, 0, 0, 0, 0),
    @('Islas Turcas y Caicos', 10, 0, 0, 9, 0, 0, 1),
    @('Cabo Verde', 10, 0, 1, 8, 0, 0, 1),
    @('Surinam', 10, 0, 6, 3, 0, 0, 1),
    @('Gambia', 9, 0, 2, 6, 0, 0, 1),
    @('Nicaragua', 9, 0, 4, 4, 0, 0, 1),
    @('Santa Sede', 8, 0, 2, 6, 0, 0, 0),
    @('Mauritania', 7, 0, 2, 4, 0, 0, 1),
    @('Sahara Occidental', 6, 0, 0, 6, 0, 0, 0),
    @('Timor Oriental', 6, 2, 1, 5, 0, 0, 0),
    @('San Bartolome', 6, 0, 4, 2, 0, 0, 0),
    @('Burundi', 5, 0, 0, 4, 0, 0, 1),
    @('Butan', 5, 0, 2, 3, 0, 0, 0),
    @('Santo Tome y Principe', 4, 0, 0, 4, 0, 0, 0),
    @('Sudan del Sur', 4, 0, 0, 4, 0, 0, 0),
    @('Bonaire, San Eustaquio y Saba', 3, 0, 0, 3, 0, 0, 0),
    @('Anguila', 3, 0, 1, 2, 0, 0, 0),
    @('Islas Virgenes Britanicas', 3, 0, 2, 1, 0, 0, 0),
    @('Papua Nueva Guinea', 2, 0, 0, 2, 0, 0, 0),
    @('San Pedro y Miquelon', 1, 0, 0, 1, 0, 0, 0),
    @('Yemen', 1, 0, 0, 1, 0, 0, 0)
)

$startRow = 4
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

$ws.Range("A1").Value = 'Datos actualizados a 14 de Abril de 2020 a las 17:52'

